$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FilePath changes from the unused PioneerNoob directory to Scene/1.xml
# (keeps its original style s="2")
$ws.Range("B2").Value = "../../NFDataCfg/Ini/Scene/1.xml"

# Rows 3-7: FilePath column was previously empty; now populated with
# Scene/2.xml .. Scene/6.xml, using a new font (family 3) / number-format (@) style.
$ws.Range("B3").Value = "../../NFDataCfg/Ini/Scene/2.xml"
$ws.Range("B3").Font.Name = "宋体"
$ws.Range("B3").Font.Size = 11
$ws.Range("B3").Font.Family = 3
$ws.Range("B3").NumberFormat = "@"

$ws.Range("B4").Value = "../../NFDataCfg/Ini/Scene/3.xml"
$ws.Range("B4").Font.Name = "宋体"
$ws.Range("B4").Font.Size = 11
$ws.Range("B4").Font.Family = 3
$ws.Range("B4").NumberFormat = "@"

$ws.Range("B5").Value = "../../NFDataCfg/Ini/Scene/4.xml"
$ws.Range("B5").Font.Name = "宋体"
$ws.Range("B5").Font.Size = 11
$ws.Range("B5").Font.Family = 3
$ws.Range("B5").NumberFormat = "@"

$ws.Range("B6").Value = "../../NFDataCfg/Ini/Scene/5.xml"
$ws.Range("B6").Font.Name = "宋体"
$ws.Range("B6").Font.Size = 11
$ws.Range("B6").Font.Family = 3
$ws.Range("B6").NumberFormat = "@"

$ws.Range("B7").Value = "../../NFDataCfg/Ini/Scene/6.xml"
$ws.Range("B7").Font.Name = "宋体"
$ws.Range("B7").Font.Size = 11
$ws.Range("B7").Font.Family = 3
$ws.Range("B7").NumberFormat = "@"

# Sheet view: scroll back to show column A (drop the old topLeftCell="C1")
# and move the active selection to B5.
$ws.Activate() | Out-Null
$ws.Range("B5").Select() | Out-Null
